$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 566, shifting existing rows 566-649 down to 567-650.
$ws.Rows.Item(566).Insert()

# Populate the newly inserted row 566 with the new data record.
$r = 566
$ws.Cells.Item($r, 1).Value = 4
$ws.Cells.Item($r, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($r, 3).Value = "Los Lagos"
$ws.Cells.Item($r, 4).Value = 45127
$ws.Cells.Item($r, 5).Value = 10
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100102
$ws.Cells.Item($r, 8).Value = "Cítricos"
$ws.Cells.Item($r, 9).Value = 100102006
$ws.Cells.Item($r, 10).Value = "Pomelo"
$ws.Cells.Item($r, 11).Value = "Start Ruby"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 100
$ws.Cells.Item($r, 14).Value = 14000
$ws.Cells.Item($r, 15).Value = 14000
$ws.Cells.Item($r, 16).Value = 14000
$ws.Cells.Item($r, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 1000
$ws.Cells.Item($r, 20).Value = 14

# Ensure the date cell keeps the date number format used by the rest of column D.
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r + 1, 4).NumberFormat
